$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Apache HTTP Server paragraph - reword the description.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Apache is one of the most widely used web servers in the world. It is open-source and supports many operating systems like Windows and Linux. It is easy to configure and reliable.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Apache is one of the most used web servers in the world. It is open source and can run on Windows and Linux. It is easy to use and very reliable.",
    2)

# ---------------------------------------------------------------------------
# 2. Nginx paragraph - reword each of the three existing runs (the
#    proofErr spellStart/spellEnd/gramStart/gramEnd markers that already
#    bracket these runs are left untouched).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    " is a very popular web server known for its high speed and performance. It is commonly used for websites with ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " is a popular web ",
    2)

$d.Content.Find.Execute(
    "a large number of",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "server",
    2)

$d.Content.Find.Execute(
    " visitors. It also works well as a load balancer.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " and it is very fast. It is mostly used for websites that have many visitors. It can also be used as a load balancer.",
    2)

# ---------------------------------------------------------------------------
# 3. Microsoft IIS paragraph.
#    a) drop the trailing space in the bold heading run.
#    b) reword the body text and split it into the same run layout used
#       elsewhere in the document: an isolated "IIS" run and an isolated
#       "in" run (mirroring the spellStart/gramStart proofed spans used for
#       the other products), surrounded by the plain narrative text.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "3. Microsoft IIS ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "3. Microsoft IIS",
    2)

$body = $d.Content
$body.Find.Execute(
    "IIS is a web server developed by Microsoft. It is mainly used on Windows servers and works well with ASP.NET applications. Many enterprise websites use IIS.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "",
    0)
$runStart = $body.Start
$newBody = "IIS is a web server made by Microsoft. It is mostly used in Windows servers. It works good with ASP.NET applications and many companies use it."
$body.Text = $newBody

# Isolate the leading "IIS" into its own run.
$iisSplit = $d.Range($runStart, $runStart + 3)
$iisSplit.Bold = 1
$iisSplit.Bold = 0

# Isolate the "in" before "Windows servers" into its own run.
$inStart = $runStart + 57
$inEnd = $runStart + 59
$inSplit = $d.Range($inStart, $inEnd)
$inSplit.Bold = 1
$inSplit.Bold = 0

# ---------------------------------------------------------------------------
# 4. LiteSpeed paragraph - reword the description.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    " is a modern and high-performance web server. It is faster than Apache and uses fewer system resources. It is commonly used for hosting WordPress websites.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " is a modern web server and has high performance. It is faster than Apache and uses less resources. It is commonly used for WordPress hosting.",
    2)
